$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.018.60'
$ws.Range('E2').Value = '  -3.44%  '
$ws.Range('D3').Value = '2.360.75'
$ws.Range('E3').Value = '  -3.88%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'500.67"
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('D6').Value = "'128.93"
$ws.Range('E6').Value = '  -3.64%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('D9').Value = '2.362.69'
$ws.Range('E9').Value = '  -3.79%  '
$ws.Range('D10').Value = "'0.0980"
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = "'4.78"
$ws.Range('E12').Value = '  +3.44%  '
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').Value = '2.777.65'
$ws.Range('E14').Value = '  -3.80%  '
$ws.Range('D15').Value = '56.022.89'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('D16').Value = "'21.37"
$ws.Range('E16').Value = '  -2.64%  '
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '2.386.86'
$ws.Range('E18').Value = '  -2.58%  '
$ws.Range('D19').Value = "'9.99"
$ws.Range('E19').Value = '  -3.46%  '
$ws.Range('D20').Value = "'4.02"
$ws.Range('E20').Value = '  -3.22%  '
$ws.Range('D21').Value = "'306.35"
$ws.Range('E21').Value = '  -2.76%  '
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = "'65.50"
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = "'0.998"
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = "'0.369"
$ws.Range('E26').Value = '  -2.80%  '
$ws.Range('D27').Value = "'0.146"
$ws.Range('E27').Value = '  -6.04%  '
$ws.Range('D28').Value = "'7.19"
$ws.Range('E28').Value = '  -5.21%  '
$ws.Range('D29').Value = "'171.41"
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').Value = '0.0₃0709'
$ws.Range('E30').Value = '  -3.53%  '
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = "'0.999"
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = "'5.73"
$ws.Range('E34').Value = '  -7.10%  '
$ws.Range('E35').Value = '  -4.95%  '
$ws.Range('D36').Value = "'17.57"
$ws.Range('E36').Value = '  -2.98%  '
$ws.Range('E37').Value = '  -6.16%  '
$ws.Range('D38').Value = "'3.73"
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('D39').Value = "'36.07"
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('D40').Value = "'0.789"
$ws.Range('E40').Value = '  -2.50%  '
$ws.Range('E41').Value = '  -5.96%  '
$ws.Range('D42').Value = "'128.78"
$ws.Range('E42').Value = '  -5.96%  '
$ws.Range('D43').Value = "'3.34"
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('D45').Value = "'0.562"
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('D46').Value = "'0.0901"
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('D47').Value = "'238.54"
$ws.Range('E47').Value = '  -7.14%  '
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('E49').Value = '  -4.06%  '
$ws.Range('D50').Value = "'17.02"
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('E51').Value = '  -0.60%  '
